# dbo_nuclide.xlsx — "helcom & de bailly handlers"
#
# Fixes the nc_name (column H) values for nuclides whose metastable-state
# marker was on the wrong side of the mass number (e.g. "mag106" -> "ag106m").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H21").Value  = "ag106m"   # 106mAg
$ws.Range("H23").Value  = "ag108m"   # 108mAg
$ws.Range("H24").Value  = "ag110m"   # 110mAg
$ws.Range("H27").Value  = "te129m"   # 129mTe
$ws.Range("H86").Value  = "ir192"    # 192Ir
$ws.Range("H91").Value  = "sn117m"   # 117mSn
$ws.Range("H92").Value  = "tl208"    # 208Tl
$ws.Range("H94").Value  = "tc99m"    # 99mTc
$ws.Range("H108").Value = "in116m"   # 116mIn
$ws.Range("H109").Value = "te123m"   # 123mTe
$ws.Range("H116").Value = "pa234m"   # 234mPa
$ws.Range("H125").Value = "ba137m"   # 137mBa

# Restore the author's on-save cursor/scroll position (row ~79, column H).
$ws.Range("H79").Select()
